# Generate Report for Handback
#
# Fills in the "Latest Target File" / "Latest Handback File" columns (E/F)
# for the two source rows on each language sheet (zh-cn, de-de), stamps the
# handback timestamp into "Latest Handback DateTime" (G), and flips the
# Status (B) from "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# BGR-packed color matching the workbook's existing hyperlink font
# (ARGB FF6495ED -> OLE color 0xED9564).
$hyperlinkColor = 15570276

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

function Handback-Sheet($sheetName, $mdUrl, $mdName, $depUrl, $depName, $xlf174Url, $xlf174Name, $xlf9e8bUrl, $xlf9e8bName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column -> "handed back" message for both source rows.
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # Row 2 (174a5250-...md source): Latest Target File / Latest Handback File
    # mirror the source markdown + handoff xlf that were already handed off
    # (columns A / C), now that they've been handed back in sync.
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl, "", "", $mdName)
    Style-AsHyperlink $ws.Range("E2")

    $ws.Hyperlinks.Add($ws.Range("F2"), $xlf174Url, "", "", $xlf174Name)
    Style-AsHyperlink $ws.Range("F2")

    $ws.Range("G2").Value = $handbackDateTime

    # Row 3 (9e8b9b78-...md source): same treatment.
    $ws.Hyperlinks.Add($ws.Range("E3"), $depUrl, "", "", $depName)
    Style-AsHyperlink $ws.Range("E3")

    $ws.Hyperlinks.Add($ws.Range("F3"), $xlf9e8bUrl, "", "", $xlf9e8bName)
    Style-AsHyperlink $ws.Range("F3")

    $ws.Range("G3").Value = $handbackDateTime
}

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/04a4ac0241766dc40877511ef8b200188fbaaeab/e2e/174a5250-aaa8-4c7e-bb4b-37d715947d46.md"
$mdName = "174a5250-aaa8-4c7e-bb4b-37d715947d46.md"
$depUrl = "https://github.com/OpenLocalizationTest/oltest/blob/04a4ac0241766dc40877511ef8b200188fbaaeab/e2e/9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.md"
$depName = "9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.md"

Handback-Sheet "zh-cn" $mdUrl $mdName $depUrl $depName `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e9dca6d6446716d89efed7f5ceb5b15eb5995130/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/174a5250-aaa8-4c7e-bb4b-37d715947d46.ab6a47d2a271508164a1a32491a6817bf22610e4.zh-cn.xlf" `
    "174a5250-aaa8-4c7e-bb4b-37d715947d46.ab6a47d2a271508164a1a32491a6817bf22610e4.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e9dca6d6446716d89efed7f5ceb5b15eb5995130/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.0499f178131fe92670f127d9da10a83ec3a86c70.zh-cn.xlf" `
    "9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.0499f178131fe92670f127d9da10a83ec3a86c70.zh-cn.xlf" `
    "2016-03-09 08:21:16"

Handback-Sheet "de-de" $mdUrl $mdName $depUrl $depName `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3fe5d31847579b497ac852962cd7f3ac95026a3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/174a5250-aaa8-4c7e-bb4b-37d715947d46.ab6a47d2a271508164a1a32491a6817bf22610e4.de-de.xlf" `
    "174a5250-aaa8-4c7e-bb4b-37d715947d46.ab6a47d2a271508164a1a32491a6817bf22610e4.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3fe5d31847579b497ac852962cd7f3ac95026a3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.0499f178131fe92670f127d9da10a83ec3a86c70.de-de.xlf" `
    "9e8b9b78-eb94-40d6-bbd0-3b23b4f28699.0499f178131fe92670f127d9da10a83ec3a86c70.de-de.xlf" `
    "2016-03-09 08:21:27"

Write-Output "Handback report generated."
